$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column B so the longer "notes" text (new row) fits
$ws.Range("B1").EntireColumn.ColumnWidth = 52.21875

# Add the new data-source row: rouses_metav_dataset (Tableau dashboard prototypes)
# Write in A, C, B, D order so new shared strings land at the same indices
# (28: name, 29: tableau url, 30: notes, 31: github url) as the authored workbook.
$ws.Range("A10").Value = "rouses_metav_dataset"
$ws.Range("C10").Value = "https://public.tableau.com/app/profile/robertrouse/viz/MetaV-app/Passages"
$ws.Range("B10").Value = "#NEED ERD# well-built dataset of the ?KJV? Bible translation"
$ws.Range("D10").Value = "https://github.com/robertrouse/theographic-bible-metadata"

# Leave the cursor where the author left it when they saved
$ws.Range("E10").Select() | Out-Null
